# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 7
    3  = 7
    4  = 5
    5  = 10
    6  = 8
    7  = 2
    8  = 8
    9  = 5
    10 = 9
    11 = 6
    12 = 5
    13 = 4
    14 = 7
    15 = 7
    16 = 8
    17 = 5
    18 = 10
    19 = 0
    20 = 2
    21 = 3
    22 = 1
    23 = 6
    24 = 8
    25 = 3
    26 = 6
    27 = 6
    28 = 5
    29 = 5
    30 = 7
    31 = 1
    32 = 6
    33 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
